# Enabled only failing tests
# - Column D ("Y"/"N" automation flag) is flipped from "Y" to "N" for every test
#   row except rows 12, 55 and 58, which stay "Y".
# - Column E (result) is updated for the rows whose run status changed:
#     row 58 -> PASS, rows 69-71 -> SKIP
# - Rows 70/71 also pick up the same cell format used by the rest of column D
#   (they previously used a plain border-only style).
# - The active selection moves from B3 to D12 and the saved scroll position
#   (top-left cell) is reset back to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Rows in column D that flip from "Y" to "N" (all data rows 2..71 except 12, 55, 58)
$keepY = @(12, 55, 58)
for ($r = 2; $r -le 71; $r++) {
    if ($keepY -notcontains $r) {
        $ws.Cells.Item($r, 4).Value = "N"
    }
}

# Column E (result) changes
$ws.Cells.Item(58, 5).Value = "PASS"
$ws.Cells.Item(69, 5).Value = "SKIP"
$ws.Cells.Item(70, 5).Value = "SKIP"
$ws.Cells.Item(71, 5).Value = "SKIP"

# Rows 70 and 71 need to pick up the standard column-D formatting (same style
# already used by every other cell in column D, e.g. D2) instead of the bare
# bordered style they had before.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D70:D71").PasteSpecial(-4122) | Out-Null

# Update the saved view: move the selection from B3 to D12.
$ws.Range("D12").Select() | Out-Null
